# Source space MEG time domain prototype
# Adds a new "lobe" column (M) to Sheet1, classifying each ROI row as
# parietal, frontal, or temporal.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M1").Value = "lobe"

# Index 0 corresponds to row 2, index 1 to row 3, etc. (rows 28-29 get no label).
$lobes = @(
    "parietal", # row 2
    "parietal", # row 3
    "frontal",  # row 4
    "frontal",  # row 5
    "parietal", # row 6
    "parietal", # row 7
    "parietal", # row 8
    "parietal", # row 9
    "parietal", # row 10
    "parietal", # row 11
    "frontal",  # row 12
    "frontal",  # row 13
    "parietal", # row 14
    "parietal", # row 15
    "parietal", # row 16
    "parietal", # row 17
    "parietal", # row 18
    "parietal", # row 19
    "parietal", # row 20
    "parietal", # row 21
    "frontal",  # row 22
    "frontal",  # row 23
    "temporal", # row 24
    "temporal", # row 25
    "temporal", # row 26
    "temporal"  # row 27
)

for ($i = 0; $i -lt $lobes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $lobes[$i]
}

$ws.Range("M28").Select()
